$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the existing "Tipo"
# header/value (column D) one column to the right, into column E,
# carrying over its formatting.
$ws.Columns.Item(4).Insert()

# Copy the formatting of the neighboring header cell (C1, bold/boxed
# header style) onto the new D1 header cell, then set its text.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "MAE"

# Populate the new MAE value for the data row.
$ws.Range("D2").Value = 0.886280647950258
